$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Berge_MDS")
$ws2 = $wb.Worksheets.Item("Reservoir_MDS")

# --- Data edits on Berge_MDS ---
# Row 2 (LE): CO2_flux_proxy_vars/subset -> rad_shortwave_down_CNR4, range -> 50
#             CH4_flux_proxy_vars/subset -> water_temp_surface
# New shared strings must be introduced in this exact order so they land at the
# same shared-string indices as the target workbook:
#   rad_shortwave_down_CNR4, air_temp, air_vpd, water_temp_surface
$ws1.Range("I2").Value = "rad_shortwave_down_CNR4"
$ws1.Range("J2").Value = "rad_shortwave_down_CNR4"

# Row 3 (H): CO2_flux_proxy_vars -> air_temp, range -> 2.5
$ws1.Range("I3").Value = "air_temp"

# Row 4 (CO2_flux): CO2_flux_proxy_vars -> air_vpd, range -> 500
$ws1.Range("I4").Value = "air_vpd"

# Back to row 2: CH4_flux_proxy_vars/subset -> water_temp_surface
$ws1.Range("L2").Value = "water_temp_surface"
$ws1.Range("M2").Value = "water_temp_surface"

# Numeric range updates
$ws1.Range("K2").Value = 50
$ws1.Range("K3").Value = 2.5
$ws1.Range("K4").Value = 500

# Row 4's CH4_flux_proxy_vars / subset / range cells are removed entirely
$ws1.Range("L4:N4").Clear()

# --- View / selection changes ---
# Berge_MDS becomes the active/selected tab, scrolled to show column K first,
# with N10 selected.
$ws1.Activate()
$ws1.Range("N10").Select()

$win = $excel.ActiveWindow
$win.ScrollColumn = 11
$win.ScrollRow = 1
